$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = -8
    7  = -4
    11 = -8
    13 = -7
    15 = -11
    17 = 5
    20 = -2
    23 = -9
    24 = 1
    28 = 0
    29 = -2
    38 = 3
    39 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
